$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "36.412.46"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.87%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.964.00"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -4.94%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "243.95"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.49%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.619"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.98%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "58.22"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -11.46%  "
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("E9").Value = "  -7.22%  "
$ws.Range("E10").Value = "  -5.94%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0846"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.62%  "
$ws.Range("E12").Value = "  -0.94%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.23"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -6.50%  "
$ws.Range("E14").Value = "  -9.70%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.251.31"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.48%  "
$ws.Range("E16").Value = "  -9.07%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.34"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -6.19%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.979.33"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.61%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "36.284.61"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.08%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "71.39"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0884"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.03%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.14"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -6.79%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "231.09"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.70%  "
$ws.Range("E24").Value = "  -0.08%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.52"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.19%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.26"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.51%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.58"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.85%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "165.94"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.25%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.91"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.31%  "
$ws.Range("E30").Value = "  -3.26%  "
$ws.Range("E31").Value = "  -3.27%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.15"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.49%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.75"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -8.64%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0641"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.56%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.33"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -7.64%  "
$ws.Range("E36").Value = "  -0.07%  "
$ws.Range("E37").Value = "  -1.76%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.96"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -7.36%  "
$ws.Range("E39").Value = "  -12.75%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.92"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.63%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0961"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -6.08%  "
$ws.Range("E42").Value = "  -5.03%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.18"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -8.33%  "
$ws.Range("E45").Value = "  -10.03%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "15.83"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -9.09%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "88.80"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -7.17%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.349.38"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.46%  "
$ws.Range("E49").Value = "  -8.79%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.82"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.46%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "45.07"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.64%  "
